$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing D-column counts (weekly refresh of source numbers) ---
$ws.Range("D2").Value = 11764
$ws.Range("D3").Value = 11509
$ws.Range("D7").Value = 11868
$ws.Range("D8").Value = 11156
$ws.Range("D12").Value = 11976
$ws.Range("D14").Value = 1321
$ws.Range("D16").Value = 1042
$ws.Range("D17").Value = 12020
$ws.Range("D19").Value = 1916
$ws.Range("D21").Value = 1500
$ws.Range("D22").Value = 12077
$ws.Range("D24").Value = 2403
$ws.Range("D26").Value = 1883
$ws.Range("D27").Value = 12119
$ws.Range("D28").Value = 9214
$ws.Range("D29").Value = 2905
$ws.Range("D31").Value = 2296
$ws.Range("D32").Value = 12159
$ws.Range("D33").Value = 8806
$ws.Range("D34").Value = 3353
$ws.Range("D36").Value = 2647
$ws.Range("D37").Value = 12197
$ws.Range("D38").Value = 8282
$ws.Range("D39").Value = 3915
$ws.Range("D41").Value = 3114
$ws.Range("D42").Value = 12232
$ws.Range("D43").Value = 7714
$ws.Range("D44").Value = 4518
$ws.Range("D46").Value = 3616
$ws.Range("D47").Value = 12257
$ws.Range("D48").Value = 7166
$ws.Range("D49").Value = 5091
$ws.Range("D51").Value = 4076
$ws.Range("D52").Value = 12285
$ws.Range("D53").Value = 6524
$ws.Range("D54").Value = 5761
$ws.Range("D56").Value = 4580
$ws.Range("D57").Value = 12308
$ws.Range("D58").Value = 5914
$ws.Range("D59").Value = 6394
$ws.Range("D61").Value = 5106
$ws.Range("D62").Value = 12319
$ws.Range("D63").Value = 5748
$ws.Range("D64").Value = 6571
$ws.Range("D66").Value = 5254
$ws.Range("D67").Value = 12341
$ws.Range("D68").Value = 5564
$ws.Range("D69").Value = 6777
$ws.Range("D71").Value = 5412
$ws.Range("D72").Value = 12357
$ws.Range("D73").Value = 5128
$ws.Range("D74").Value = 7229
$ws.Range("D76").Value = 5813
$ws.Range("D77").Value = 12380
$ws.Range("D78").Value = 4651
$ws.Range("D79").Value = 7729
$ws.Range("D81").Value = 6257
$ws.Range("D82").Value = 12398
$ws.Range("D84").Value = 8266
$ws.Range("D86").Value = 6758
$ws.Range("D87").Value = 12425
$ws.Range("D88").Value = 3726
$ws.Range("D89").Value = 8699
$ws.Range("D90").Value = 1538
$ws.Range("D91").Value = 7161

# --- Append new week (202506 / 2025-02-09) rows 92-96 ---
# Copy the date-formatted style from the last existing date cell (B91) so the
# new date cells reuse the same number-format style index instead of creating a new one.
$ws.Range("B91").Copy()
$ws.Range("B92").PasteSpecial(-4122)
$ws.Cells.Item(92, 1).Value = 202506
$ws.Cells.Item(92, 2).Value = 45697
$ws.Cells.Item(92, 3).Value = "farms_total_count"
$ws.Cells.Item(92, 4).Value = 12441

$ws.Range("B91").Copy()
$ws.Range("B93").PasteSpecial(-4122)
$ws.Cells.Item(93, 1).Value = 202506
$ws.Cells.Item(93, 2).Value = 45697
$ws.Cells.Item(93, 3).Value = "farms_to_examine_count"
$ws.Cells.Item(93, 4).Value = 3388

$ws.Range("B91").Copy()
$ws.Range("B94").PasteSpecial(-4122)
$ws.Cells.Item(94, 1).Value = 202506
$ws.Cells.Item(94, 2).Value = 45697
$ws.Cells.Item(94, 3).Value = "farms_examined_count"
$ws.Cells.Item(94, 4).Value = 9053

$ws.Range("B91").Copy()
$ws.Range("B95").PasteSpecial(-4122)
$ws.Cells.Item(95, 1).Value = 202506
$ws.Cells.Item(95, 2).Value = 45697
$ws.Cells.Item(95, 3).Value = "farms_examined_positive_count"
$ws.Cells.Item(95, 4).Value = 1534

$ws.Range("B91").Copy()
$ws.Range("B96").PasteSpecial(-4122)
$ws.Cells.Item(96, 1).Value = 202506
$ws.Cells.Item(96, 2).Value = 45697
$ws.Cells.Item(96, 3).Value = "farms_examined_negative_count"
$ws.Cells.Item(96, 4).Value = 7519

# --- Clear clipboard marching ants left over from copy operations ---
$excel.CutCopyMode = 0

# --- Selection state matches the saved workbook (whole-sheet selection) ---
$ws.Cells.Select() | Out-Null

Write-Host "Applied Moderhinke Dashboard weekly update"
